$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    288504.5508113568,
    37609.82862424095,
    75565.8528822182,
    78949.95940588829,
    61159.02891982916,
    209731.5255150525,
    69626.30223623497,
    46469.70857545878,
    35696.14385609442,
    63123.84478749041,
    53464.91023357871,
    81342.73572248663,
    81739.52457452597,
    19276.85656763474,
    117630.9536402359,
    99251.3719825726,
    15979.75015829244,
    82927.09302641748,
    45393.46432031956,
    64737.73085822884,
    49632.69356745129,
    41955.87687536592,
    86807.04049347091,
    88327.06810272801,
    13387.2039112985,
    52434.92956434948,
    56425.53542155959,
    75828.26811740537,
    78731.49033877561,
    188366.5224977621,
    68457.13245405869,
    64582.49572546703,
    103908.6745720124,
    113678.5499176987,
    94354.3146973332,
    78231.57813929534,
    85374.78615633408,
    22826.92162173722,
    231684.2522995113,
    92134.81420432765,
    49637.13083274924,
    62466.95265461144,
    37925.66687774626,
    127668.1739258605,
    10681.45704567587,
    52800.57197574813,
    43846.64443427842,
    54480.28487128675,
    105453.8803313497,
    41732.84436106237,
    59185.56625616069,
    105935.0833999371,
    32639.47633890131,
    35578.43318048371,
    28872.66159863357,
    193509.3354416623,
    44183.33248351423,
    20126.40780746128,
    67991.28610770732,
    108611.6782937965,
    215771.5200893505,
    27211.84075890353,
    221896.2272498873,
    116860.7111063749,
    87068.56717412683,
    77127.90197223205,
    51208.20564629497,
    71891.96432015175,
    96523.56822806783,
    28014.64430274286,
    126663.9479380092,
    29510.65657633863,
    51809.78794763309,
    216422.2744131096,
    29344.22800173918,
    86325.89188858398,
    45211.55894453867,
    70353.6320838845,
    96143.76320364367,
    43331.71188431633,
    145844.8463954173,
    45967.54278047707,
    30561.36720004267,
    20249.21186314673,
    188687.3793641133,
    92686.66370008531,
    86460.60662827775,
    30336.23770554819,
    42396.88328025799,
    142429.2510757815,
    62840.98370790164,
    36088.37372944201,
    78887.39844741176,
    61317.30057170129,
    367767.7875245648,
    126952.1394098008,
    57549.78057554606,
    63200.69823052623,
    4925.772066883168,
    61362.41111483047
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $values[$i]
}

Write-Output "Updated $($values.Length) residual values in column A (A2:A101)"
